$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 3-5 (keep header row 1 and first data row 2)
$ws.Range("A3:D5").EntireRow.Delete() | Out-Null

# Update the remaining data row with new values
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 10000
$ws.Range("D2").Value = 20

# Update the active selection to match the target state
$ws.Range("F6").Select() | Out-Null
